# Apply the symbol-list update: numeric/percentage cells in columns D and E get a
# leading apostrophe so Excel stores them as text (matching the workbook's existing
# inline-string convention) instead of auto-converting to numbers/percent values.
# Columns B (coin name) and C (link) are plain text already, so no prefix is needed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'297.48"
$ws.Range("E2").Value = "'2.57%"
$ws.Range("D3").Value = "'41.48"
$ws.Range("E3").Value = "'2.63%"
$ws.Range("D4").Value = "'5.021"
$ws.Range("E4").Value = "'-0.44%"
$ws.Range("E5").Value = "'3.28%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.376"
$ws.Range("E6").Value = "'2.31%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.587"
$ws.Range("E7").Value = "'1.37%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9294"
$ws.Range("E8").Value = "'1.30%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.423"
$ws.Range("E9").Value = "'1.09%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1198"
$ws.Range("E10").Value = "'2.42%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1843"
$ws.Range("E11").Value = "'7.34%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08909"
$ws.Range("E12").Value = "'3.63%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.04037"
$ws.Range("E13").Value = "'-3.15%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1054"
$ws.Range("E14").Value = "'0.11%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001276"
$ws.Range("E15").Value = "'-0.22%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005839"
$ws.Range("E16").Value = "'1.19%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.342"
$ws.Range("E17").Value = "'-1.56%"
$ws.Range("D18").Value = "'0.3311"
$ws.Range("E18").Value = "'1.08%"
$ws.Range("D19").Value = "'7.990"
$ws.Range("E19").Value = "'1.89%"
$ws.Range("D20").Value = "'0.1418"
$ws.Range("E20").Value = "'5.04%"
$ws.Range("E21").Value = "'4.05%"
$ws.Range("D22").Value = "'0.04054"
$ws.Range("E22").Value = "'4.97%"
$ws.Range("E23").Value = "'-0.15%"
$ws.Range("D24").Value = "'0.004235"
$ws.Range("E24").Value = "'9.94%"
$ws.Range("D25").Value = "'0.0001230"
$ws.Range("E25").Value = "'-3.98%"
$ws.Range("E26").Value = "'0.04%"
$ws.Range("D38").Value = "'0.02417"
$ws.Range("E38").Value = "'4.46%"
$ws.Range("D39").Value = "'0.05216"
$ws.Range("E39").Value = "'4.89%"
$ws.Range("D40").Value = "'0.006558"
$ws.Range("E40").Value = "'-2.43%"
$ws.Range("D41").Value = "'0.007793"
$ws.Range("E41").Value = "'1.54%"
$ws.Range("D42").Value = "'0.1331"
$ws.Range("E42").Value = "'4.21%"
$ws.Range("D43").Value = "'0.007637"
$ws.Range("E43").Value = "'3.67%"
$ws.Range("D44").Value = "'0.007826"
$ws.Range("E44").Value = "'10.76%"
$ws.Range("D45").Value = "'0.3219"
$ws.Range("E45").Value = "'11.44%"
$ws.Range("D46").Value = "'0.00006689"
$ws.Range("E46").Value = "'4.05%"
$ws.Range("E47").Value = "'-0.05%"
$ws.Range("E48").Value = "'-0.07%"
$ws.Range("E49").Value = "'442.00%"
$ws.Range("E50").Value = "'-0.05%"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'-0.05%"
